$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.495.36"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").Value = "1.642.83"
$ws.Range("E3").Value = "  +0.97%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("D6").Value = "'303.36"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").Value = "'0.3829"
$ws.Range("E7").Value = "  +1.66%  "

$ws.Range("D8").Value = "'51.99"
$ws.Range("E8").Value = "  -0.36%  "

$ws.Range("D9").Value = "'0.3606"
$ws.Range("E9").Value = "  -0.40%  "

$ws.Range("D10").Value = "'0.08261"
$ws.Range("E10").Value = "  +2.24%  "

$ws.Range("D11").Value = "'1.232"
$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("D13").Value = "'22.55"
$ws.Range("E13").Value = "  -0.60%  "

$ws.Range("D14").Value = "'6.463"
$ws.Range("E14").Value = "  -1.48%  "

$ws.Range("D15").Value = "'7.357"
$ws.Range("E15").Value = "  +1.86%  "

$ws.Range("E16").Value = "  -0.52%  "

$ws.Range("D17").Value = "1.637.56"
$ws.Range("E17").Value = "  +0.77%  "

$ws.Range("D18").Value = "'95.19"
$ws.Range("E18").Value = "  +1.88%  "

$ws.Range("D19").Value = "'0.06970"
$ws.Range("E19").Value = "  +0.59%  "

$ws.Range("D20").Value = "'6.598"
$ws.Range("E20").Value = "  +2.13%  "

$ws.Range("D21").Value = "'17.51"
$ws.Range("E21").Value = "  -2.45%  "

$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").Value = "'12.53"
$ws.Range("E23").Value = "  -1.63%  "

$ws.Range("D24").Value = "23.477.14"
$ws.Range("E24").Value = "  +0.42%  "

$ws.Range("D25").Value = "'2.535"
$ws.Range("E25").Value = "  +4.76%  "

$ws.Range("D26").Value = "'3.078"
$ws.Range("E26").Value = "  -4.73%  "

$ws.Range("D27").Value = "'21.17"
$ws.Range("E27").Value = "  +0.17%  "

$ws.Range("D28").Value = "'151.93"
$ws.Range("E28").Value = "  +2.00%  "

$ws.Range("D29").Value = "'5.278"
$ws.Range("E29").Value = "  -0.51%  "

$ws.Range("D30").Value = "'133.46"
$ws.Range("E30").Value = "  -1.04%  "

$ws.Range("D31").Value = "1.820.51"
$ws.Range("E31").Value = "  +0.77%  "

$ws.Range("D32").Value = "'1.087"
$ws.Range("E32").Value = "  +14.45%  "

$ws.Range("D33").Value = "'2.154"
$ws.Range("E33").Value = "  -6.45%  "

$ws.Range("D34").Value = "'6.536"
$ws.Range("E34").Value = "  -3.75%  "

$ws.Range("D35").Value = "'11.50"
$ws.Range("E35").Value = "  +5.38%  "

$ws.Range("E36").Value = "  -2.04%  "

$ws.Range("D37").Value = "'0.2515"
$ws.Range("E37").Value = "  -0.69%  "

$ws.Range("D38").Value = "'0.08794"
$ws.Range("E38").Value = "  -0.29%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.07039"
$ws.Range("E39").Value = "  -1.13%  "

$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "'5.972"
$ws.Range("E40").Value = "  -2.71%  "

$ws.Range("D41").Value = "'0.7039"
$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("D42").Value = "'1.348"
$ws.Range("E42").Value = "  -0.94%  "

$ws.Range("D43").Value = "'12.28"
$ws.Range("E43").Value = "  -0.63%  "

$ws.Range("D44").Value = "'15.61"
$ws.Range("E44").Value = "  -3.70%  "

$ws.Range("D45").Value = "'0.6524"
$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").Value = "'2.295"
$ws.Range("E47").Value = "  -1.13%  "

$ws.Range("D48").Value = "'3.961"
$ws.Range("E48").Value = "  -0.51%  "

$ws.Range("D49").Value = "'0.07985"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").Value = "'128.99"
$ws.Range("E50").Value = "  +2.26%  "

$ws.Range("E51").Value = "  -1.14%  "
